$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text labels (translated from Danish to English)
$ws.Range("A1").Value = "Individuals use of cultural activities (year) within the past three months by time, area and cultural activities"
$ws.Range("A2").Value = "Units: Per cent"
$ws.Range("C3").Value = "Have been to the cinema"
$ws.Range("D3").Value = "Have been to concert"
$ws.Range("E3").Value = "Have seen visual arts (on pupose)"
$ws.Range("F3").Value = "Have read or listened to fiction"
$ws.Range("G3").Value = "Have read or listened to non-fiction"
$ws.Range("H3").Value = "Have visited the library (physical visit)"
$ws.Range("I3").Value = "Have used the librarys digital services"
$ws.Range("J3").Value = "Visited a museum etc."
$ws.Range("K3").Value = "Have watched performing arts in theater, opera, festivals or in public spaces"
$ws.Range("B5").Value = "Province Byen København"
$ws.Range("B6").Value = "Copenhagen"
$ws.Range("B8").Value = "Province Københavns omegn"
$ws.Range("B11").Value = "Province Nordsjælland"
$ws.Range("B16").Value = "Province Bornholm"
$ws.Range("B18").Value = "Province Østsjælland"
$ws.Range("B22").Value = "Province Vest- og Sydsjælland"
$ws.Range("B29").Value = "Province Fyn"
$ws.Range("B34").Value = "Province Sydjylland"
$ws.Range("B44").Value = "Province Østjylland"
$ws.Range("B52").Value = "Province Vestjylland"
$ws.Range("B59").Value = "Province Nordjylland"

# Cells whose text starts with a literal apostrophe need special handling,
# since assigning to .Value directly causes Excel to treat a leading
# apostrophe as a text-prefix marker and strip it. Build the string via a
# formula (CHAR(39) for the apostrophe) and then convert the formula result
# to a plain value with PasteSpecial so the cell ends up as static text.
$ws.Range("A65").Formula = "=CHAR(39)&`"Have seen movies and series' covers watching at home or outside and includes 'Have been in the cinema'. 'Have listened to music' is music at home or out of home and includes 'Have been to a concert'. 'Have watched or listened to a sporting event' is sport watched at home eg on television and/or in the city, and includes the activity 'Been to a sport event as a spectator'. Discretioniced municipalities and Bornholm have few observations `""
$ws.Range("A65").Copy() | Out-Null
$ws.Range("A65").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
$excel.CutCopyMode = 0

# Update column widths (B..J) - set ColumnWidth as close as possible to target
$ws.Columns.Item(2).ColumnWidth = 27.717291666666668
$ws.Columns.Item(3).ColumnWidth = 24.147291666666668
$ws.Columns.Item(4).ColumnWidth = 20.867291666666667
$ws.Columns.Item(5).ColumnWidth = 31.717291666666664
$ws.Columns.Item(6).ColumnWidth = 29.29729166666667
$ws.Columns.Item(7).ColumnWidth = 33.43729166666667
$ws.Columns.Item(8).ColumnWidth = 35.71729166666667
$ws.Columns.Item(9).ColumnWidth = 35.297291666666666
$ws.Columns.Item(10).ColumnWidth = 21.717291666666668
